$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.757063508033752
$ws.Range("B1").Value = 1.755508542060852
$ws.Range("C1").Value = 2.181628704071045
$ws.Range("D1").Value = 2.003978252410889
$ws.Range("E1").Value = 2.976402282714844
